# nst-est2019-01.xlsx : "update to ER data wrangling"
#
# The state-name column (A10:A60) used a "hidden leading dot" trick: each
# shared string was a rich-text run of an invisible "." (white font) followed
# by the visible state name run. This edit removes that trick: the cells now
# hold the plain state name as a simple string, with normal (visible) text
# color instead of white/invisible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$states = @(
    "Alabama", "Alaska", "Arizona", "Arkansas", "California", "Colorado",
    "Connecticut", "Delaware", "District of Columbia", "Florida", "Georgia",
    "Hawaii", "Idaho", "Illinois", "Indiana", "Iowa", "Kansas", "Kentucky",
    "Louisiana", "Maine", "Maryland", "Massachusetts", "Michigan", "Minnesota",
    "Mississippi", "Missouri", "Montana", "Nebraska", "Nevada",
    "New Hampshire", "New Jersey", "New Mexico", "New York", "North Carolina",
    "North Dakota", "Ohio", "Oklahoma", "Oregon", "Pennsylvania",
    "Rhode Island", "South Carolina", "South Dakota", "Tennessee", "Texas",
    "Utah", "Vermont", "Virginia", "Washington", "West Virginia", "Wisconsin",
    "Wyoming"
)

$startRow = 10
for ($i = 0; $i -lt $states.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $states[$i]
    $cell.Font.Name = "MS sans serif"
    $cell.Font.Size = 10
    $cell.Font.Color = 0
}

# Restore the cursor/selection to where the saved file shows it (A7).
$ws.Range("A7").Select()
